$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the second average column
$ws.Range("M2").Value = "Avarage"

# Row 4: plain (non-shared) formulas
$ws.Range("L4").Formula = "=AVERAGE(B4,D4,F4,H4,J4)"
$ws.Range("M4").Formula = "=AVERAGE(C4,E4,G4,I4,K4)"

# Rows 5-8: same formulas (will be authored as shared formulas on save)
$ws.Range("L5").Formula = "=AVERAGE(B5,D5,F5,H5,J5)"
$ws.Range("M5").Formula = "=AVERAGE(C5,E5,G5,I5,K5)"

$ws.Range("L6").Formula = "=AVERAGE(B6,D6,F6,H6,J6)"
$ws.Range("M6").Formula = "=AVERAGE(C6,E6,G6,I6,K6)"

$ws.Range("L7").Formula = "=AVERAGE(B7,D7,F7,H7,J7)"
$ws.Range("M7").Formula = "=AVERAGE(C7,E7,G7,I7,K7)"

$ws.Range("L8").Formula = "=AVERAGE(B8,D8,F8,H8,J8)"
$ws.Range("M8").Formula = "=AVERAGE(C8,E8,G8,I8,K8)"

# Selection, as recorded by the UI after entering the formulas
$ws.Range("M4:M8").Select()
